# "add first class soul" - adds a new enemy row (id=3) with a "first class soul"
# (soulid 5000) and resets several rows-2/3 stat columns, per the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- Update existing row 2 (enemy id=1) ----
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("L2").Value = 1000

# ---- Update existing row 3 (enemy id=2) ----
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("I3").Value = 1
$ws.Range("L3").Value = 1001

# ---- Add new row 4 (enemy id=3) ----
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "ENEMY_NAME_3"
$ws.Range("C4").Value = "ENEMY_DESC_3"
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 2400
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 5000
$ws.Range("M4").Value = 1
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 1
$ws.Range("R4").Value = 1
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 2
$ws.Range("W4").Value = 1
$ws.Range("X4").Value = 1
$ws.Range("Y4").Value = 3
$ws.Range("Z4").Value = 1
$ws.Range("AA4").Value = 1

# ---- Fix up the view: drop the scrolled topLeftCell and move the selection ----
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("J14").Select()
